$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'ECs'
$ws.Cells.Item(2, 2).Value = 'Bmp7'
$ws.Cells.Item(2, 3).Value = 'Acvr1'
$ws.Cells.Item(2, 4).Value = 'ECs'
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.03443933333333333
$ws.Cells.Item(2, 8).Value = 0.103318
$ws.Cells.Item(2, 9).Value = 0.05823261822459219
$ws.Cells.Item(2, 10).Value = 0.0582326182245922
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 8.165540666666667
$ws.Cells.Item(2, 14).Value = 24.496622
$ws.Cells.Item(2, 15).Value = 0.1715865889461355
$ws.Cells.Item(2, 16).Value = 0.1715865889461355
$ws.Cells.Item(2, 17).Value = 0.2812157768662222
$ws.Cells.Item(2, 18).Value = 2.530941991796
$ws.Cells.Item(2, 19).Value = 0.009991936326560338
$ws.Cells.Item(2, 20).Value = 0.009991936326560341

$ws.Cells.Item(3, 1).Value = 'ECs'
$ws.Cells.Item(3, 2).Value = 'Bmp7'
$ws.Cells.Item(3, 3).Value = 'Acvr1'
$ws.Cells.Item(3, 4).Value = 'FAPs'
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.03443933333333333
$ws.Cells.Item(3, 8).Value = 0.103318
$ws.Cells.Item(3, 9).Value = 0.05823261822459219
$ws.Cells.Item(3, 10).Value = 0.0582326182245922
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 22.011801
$ws.Cells.Item(3, 14).Value = 66.035403
$ws.Cells.Item(3, 15).Value = 0.4625449807101323
$ws.Cells.Item(3, 16).Value = 0.4625449807101323
$ws.Cells.Item(3, 17).Value = 0.758071751906
$ws.Cells.Item(3, 18).Value = 6.822645767154
$ws.Cells.Item(3, 19).Value = 0.02693520527339449
$ws.Cells.Item(3, 20).Value = 0.0269352052733945

$ws.Cells.Item(4, 1).Value = 'ECs'
$ws.Cells.Item(4, 2).Value = 'Bmp7'
$ws.Cells.Item(4, 3).Value = 'Acvr1'
$ws.Cells.Item(4, 4).Value = 'MuSCs'
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.03443933333333333
$ws.Cells.Item(4, 8).Value = 0.103318
$ws.Cells.Item(4, 9).Value = 0.05823261822459219
$ws.Cells.Item(4, 10).Value = 0.0582326182245922
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 13.51552533333333
$ws.Cells.Item(4, 14).Value = 40.546576
$ws.Cells.Item(4, 15).Value = 0.2840084918355372
$ws.Cells.Item(4, 16).Value = 0.2840084918355373
$ws.Cells.Item(4, 17).Value = 0.4654656821297777
$ws.Cells.Item(4, 18).Value = 4.189191139168
$ws.Cells.Item(4, 19).Value = 0.01653855807760105
$ws.Cells.Item(4, 20).Value = 0.01653855807760105

$ws.Cells.Item(5, 1).Value = 'ECs'
$ws.Cells.Item(5, 2).Value = 'Bmp7'
$ws.Cells.Item(5, 3).Value = 'Acvr1'
$ws.Cells.Item(5, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.03443933333333333
$ws.Cells.Item(5, 8).Value = 0.103318
$ws.Cells.Item(5, 9).Value = 0.05823261822459219
$ws.Cells.Item(5, 10).Value = 0.0582326182245922
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 3.895588
$ws.Cells.Item(5, 14).Value = 11.686764
$ws.Cells.Item(5, 15).Value = 0.08185993850819488
$ws.Cells.Item(5, 16).Value = 0.0818599385081949
$ws.Cells.Item(5, 17).Value = 0.1341614536613333
$ws.Cells.Item(5, 18).Value = 1.207453082952
$ws.Cells.Item(5, 19).Value = 0.004766918547036305
$ws.Cells.Item(5, 20).Value = 0.004766918547036307

$ws.Cells.Item(6, 1).Value = 'FAPs'
$ws.Cells.Item(6, 2).Value = 'Bmp7'
$ws.Cells.Item(6, 3).Value = 'Acvr1'
$ws.Cells.Item(6, 4).Value = 'ECs'
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.51625
$ws.Cells.Item(6, 8).Value = 1.54875
$ws.Cells.Item(6, 9).Value = 0.8729143757654733
$ws.Cells.Item(6, 10).Value = 0.8729143757654734
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 8.165540666666667
$ws.Cells.Item(6, 14).Value = 24.496622
$ws.Cells.Item(6, 15).Value = 0.1715865889461355
$ws.Cells.Item(6, 16).Value = 0.1715865889461355
$ws.Cells.Item(6, 17).Value = 4.215460369166666
$ws.Cells.Item(6, 18).Value = 37.93914332250001
$ws.Cells.Item(6, 19).Value = 0.1497804001796427
$ws.Cells.Item(6, 20).Value = 0.1497804001796428

$ws.Cells.Item(7, 1).Value = 'FAPs'
$ws.Cells.Item(7, 2).Value = 'Bmp7'
$ws.Cells.Item(7, 3).Value = 'Acvr1'
$ws.Cells.Item(7, 4).Value = 'FAPs'
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.51625
$ws.Cells.Item(7, 8).Value = 1.54875
$ws.Cells.Item(7, 9).Value = 0.8729143757654733
$ws.Cells.Item(7, 10).Value = 0.8729143757654734
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 22.011801
$ws.Cells.Item(7, 14).Value = 66.035403
$ws.Cells.Item(7, 15).Value = 0.4625449807101323
$ws.Cells.Item(7, 16).Value = 0.4625449807101323
$ws.Cells.Item(7, 17).Value = 11.36359226625
$ws.Cells.Item(7, 18).Value = 102.27233039625
$ws.Cells.Item(7, 19).Value = 0.403762163100038
$ws.Cells.Item(7, 20).Value = 0.4037621631000381

$ws.Cells.Item(8, 1).Value = 'FAPs'
$ws.Cells.Item(8, 2).Value = 'Bmp7'
$ws.Cells.Item(8, 3).Value = 'Acvr1'
$ws.Cells.Item(8, 4).Value = 'MuSCs'
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.51625
$ws.Cells.Item(8, 8).Value = 1.54875
$ws.Cells.Item(8, 9).Value = 0.8729143757654733
$ws.Cells.Item(8, 10).Value = 0.8729143757654734
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 13.51552533333333
$ws.Cells.Item(8, 14).Value = 40.546576
$ws.Cells.Item(8, 15).Value = 0.2840084918355372
$ws.Cells.Item(8, 16).Value = 0.2840084918355373
$ws.Cells.Item(8, 17).Value = 6.977389953333334
$ws.Cells.Item(8, 18).Value = 62.79650958000001
$ws.Cells.Item(8, 19).Value = 0.2479150953627115
$ws.Cells.Item(8, 20).Value = 0.2479150953627116

$ws.Cells.Item(9, 1).Value = 'FAPs'
$ws.Cells.Item(9, 2).Value = 'Bmp7'
$ws.Cells.Item(9, 3).Value = 'Acvr1'
$ws.Cells.Item(9, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.51625
$ws.Cells.Item(9, 8).Value = 1.54875
$ws.Cells.Item(9, 9).Value = 0.8729143757654733
$ws.Cells.Item(9, 10).Value = 0.8729143757654734
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.895588
$ws.Cells.Item(9, 14).Value = 11.686764
$ws.Cells.Item(9, 15).Value = 0.08185993850819488
$ws.Cells.Item(9, 16).Value = 0.0818599385081949
$ws.Cells.Item(9, 17).Value = 2.011097305
$ws.Cells.Item(9, 18).Value = 18.099875745
$ws.Cells.Item(9, 19).Value = 0.07145671712308097
$ws.Cells.Item(9, 20).Value = 0.07145671712308098

$ws.Cells.Item(10, 1).Value = 'MuSCs'
$ws.Cells.Item(10, 2).Value = 'Bmp7'
$ws.Cells.Item(10, 3).Value = 'Acvr1'
$ws.Cells.Item(10, 4).Value = 'ECs'
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.04072033333333334
$ws.Cells.Item(10, 8).Value = 0.122161
$ws.Cells.Item(10, 9).Value = 0.06885300600993445
$ws.Cells.Item(10, 10).Value = 0.06885300600993445
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 8.165540666666667
$ws.Cells.Item(10, 14).Value = 24.496622
$ws.Cells.Item(10, 15).Value = 0.1715865889461355
$ws.Cells.Item(10, 16).Value = 0.1715865889461355
$ws.Cells.Item(10, 17).Value = 0.3325035377935556
$ws.Cells.Item(10, 18).Value = 2.992531840142
$ws.Cells.Item(10, 19).Value = 0.01181425243993242
$ws.Cells.Item(10, 20).Value = 0.01181425243993242

$ws.Cells.Item(11, 1).Value = 'MuSCs'
$ws.Cells.Item(11, 2).Value = 'Bmp7'
$ws.Cells.Item(11, 3).Value = 'Acvr1'
$ws.Cells.Item(11, 4).Value = 'FAPs'
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.04072033333333334
$ws.Cells.Item(11, 8).Value = 0.122161
$ws.Cells.Item(11, 9).Value = 0.06885300600993445
$ws.Cells.Item(11, 10).Value = 0.06885300600993445
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 22.011801
$ws.Cells.Item(11, 14).Value = 66.035403
$ws.Cells.Item(11, 15).Value = 0.4625449807101323
$ws.Cells.Item(11, 16).Value = 0.4625449807101323
$ws.Cells.Item(11, 17).Value = 0.8963278739870002
$ws.Cells.Item(11, 18).Value = 8.066950865883001
$ws.Cells.Item(11, 19).Value = 0.03184761233669976
$ws.Cells.Item(11, 20).Value = 0.03184761233669976

$ws.Cells.Item(12, 1).Value = 'MuSCs'
$ws.Cells.Item(12, 2).Value = 'Bmp7'
$ws.Cells.Item(12, 3).Value = 'Acvr1'
$ws.Cells.Item(12, 4).Value = 'MuSCs'
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.04072033333333334
$ws.Cells.Item(12, 8).Value = 0.122161
$ws.Cells.Item(12, 9).Value = 0.06885300600993445
$ws.Cells.Item(12, 10).Value = 0.06885300600993445
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 13.51552533333333
$ws.Cells.Item(12, 14).Value = 40.546576
$ws.Cells.Item(12, 15).Value = 0.2840084918355372
$ws.Cells.Item(12, 16).Value = 0.2840084918355373
$ws.Cells.Item(12, 17).Value = 0.5503566967484446
$ws.Cells.Item(12, 18).Value = 4.953210270736
$ws.Cells.Item(12, 19).Value = 0.01955483839522466
$ws.Cells.Item(12, 20).Value = 0.01955483839522467

$ws.Cells.Item(13, 1).Value = 'MuSCs'
$ws.Cells.Item(13, 2).Value = 'Bmp7'
$ws.Cells.Item(13, 3).Value = 'Acvr1'
$ws.Cells.Item(13, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.04072033333333334
$ws.Cells.Item(13, 8).Value = 0.122161
$ws.Cells.Item(13, 9).Value = 0.06885300600993445
$ws.Cells.Item(13, 10).Value = 0.06885300600993445
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 3.895588
$ws.Cells.Item(13, 14).Value = 11.686764
$ws.Cells.Item(13, 15).Value = 0.08185993850819488
$ws.Cells.Item(13, 16).Value = 0.0818599385081949
$ws.Cells.Item(13, 17).Value = 0.1586296418893333
$ws.Cells.Item(13, 18).Value = 1.427666777004
$ws.Cells.Item(13, 19).Value = 0.005636302838077607
$ws.Cells.Item(13, 20).Value = 0.005636302838077607
